# Auto-generated edit script: updates Leve price/profit columns (H-N)
# across all 8 profession sheets based on refreshed market-board data.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3949.8
$ws.Range("J17").Value = 4687.25
$ws.Range("L17").Value = 14061.75
$ws.Range("N17").Value = -14397.75
$ws.Range("H51").Value = 12303.263
$ws.Range("I51").Value = 16359.2
$ws.Range("J51").Value = 10854.714
$ws.Range("K51").Value = 16359.2
$ws.Range("L51").Value = 10854.714
$ws.Range("M51").Value = -15875.2
$ws.Range("N51").Value = -11822.714
$ws.Range("H70").Value = 3186.0417
$ws.Range("I70").Value = 1698
$ws.Range("J70").Value = 3930.0625
$ws.Range("K70").Value = 5094
$ws.Range("L70").Value = 11790.1875
$ws.Range("M70").Value = -4824
$ws.Range("N70").Value = -12330.1875
$ws.Range("H73").Value = 3186.0417
$ws.Range("I73").Value = 1698
$ws.Range("J73").Value = 3930.0625
$ws.Range("K73").Value = 5094
$ws.Range("L73").Value = 11790.1875
$ws.Range("M73").Value = -4158
$ws.Range("N73").Value = -13662.1875
$ws.Range("H113").Value = 3079.3914
$ws.Range("I113").Value = 2845.6
$ws.Range("J113").Value = 3144.3333
$ws.Range("K113").Value = 2845.6
$ws.Range("L113").Value = 3144.3333
$ws.Range("M113").Value = 408.4000000000001
$ws.Range("N113").Value = -9652.3333
$ws.Range("H138").Value = 2427.9844
$ws.Range("J138").Value = 2892.5454
$ws.Range("L138").Value = 8677.636200000001
$ws.Range("N138").Value = -18957.6362

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8625.362999999999
$ws.Range("J74").Value = 12594.5
$ws.Range("L74").Value = 12594.5
$ws.Range("N74").Value = -14342.5
$ws.Range("H77").Value = 8625.362999999999
$ws.Range("J77").Value = 12594.5
$ws.Range("L77").Value = 62972.5
$ws.Range("N77").Value = -71708.5
$ws.Range("H97").Value = 1030.1052
$ws.Range("I97").Value = 984.55554
$ws.Range("K97").Value = 984.55554
$ws.Range("M97").Value = -488.55554
$ws.Range("H102").Value = 1069.5714
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2053.311
$ws.Range("I134").Value = 1748.1464
$ws.Range("K134").Value = 5244.439200000001
$ws.Range("M134").Value = -2709.439200000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2392.0667
$ws.Range("I22").Value = 2652.818
$ws.Range("J22").Value = 1675
$ws.Range("K22").Value = 2652.818
$ws.Range("L22").Value = 1675
$ws.Range("M22").Value = -2302.818
$ws.Range("N22").Value = -2375
$ws.Range("H50").Value = 32633
$ws.Range("J50").Value = 32633
$ws.Range("L50").Value = 32633
$ws.Range("N50").Value = -33883
$ws.Range("H51").Value = 18412.572
$ws.Range("J51").Value = 18499.666
$ws.Range("L51").Value = 18499.666
$ws.Range("N51").Value = -19971.666
$ws.Range("H58").Value = 3578.6428
$ws.Range("I58").Value = 3131.1667
$ws.Range("K58").Value = 3131.1667
$ws.Range("M58").Value = -2928.1667
$ws.Range("H59").Value = 106999.75
$ws.Range("I59").Value = 90000
$ws.Range("J59").Value = 112666.336
$ws.Range("K59").Value = 90000
$ws.Range("L59").Value = 112666.336
$ws.Range("M59").Value = -88855
$ws.Range("N59").Value = -114956.336
$ws.Range("H60").Value = 27764.412
$ws.Range("I60").Value = 10015.8
$ws.Range("J60").Value = 35159.668
$ws.Range("K60").Value = 10015.8
$ws.Range("L60").Value = 35159.668
$ws.Range("M60").Value = -9504.799999999999
$ws.Range("N60").Value = -36181.668
$ws.Range("H61").Value = 18412.572
$ws.Range("J61").Value = 18499.666
$ws.Range("L61").Value = 18499.666
$ws.Range("N61").Value = -19195.666
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H75").Value = 60257
$ws.Range("J75").Value = 60257
$ws.Range("L75").Value = 60257
$ws.Range("N75").Value = -62253
$ws.Range("H78").Value = 60257
$ws.Range("J78").Value = 60257
$ws.Range("L78").Value = 180771
$ws.Range("N78").Value = -190755
$ws.Range("H94").Value = 1862.4
$ws.Range("J94").Value = 2114.125
$ws.Range("L94").Value = 2114.125
$ws.Range("N94").Value = -3016.125
$ws.Range("H132").Value = 4878.3076
$ws.Range("I132").Value = 2929.9062
$ws.Range("K132").Value = 8789.7186
$ws.Range("M132").Value = -6259.7186
$ws.Range("H135").Value = 123999
$ws.Range("J135").Value = 123999
$ws.Range("L135").Value = 123999
$ws.Range("N135").Value = -134139
$ws.Range("H136").Value = 3578.6428
$ws.Range("I136").Value = 3131.1667
$ws.Range("K136").Value = 9393.500100000001
$ws.Range("M136").Value = -6843.500100000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 15912.375
$ws.Range("J39").Value = 17471.285
$ws.Range("L39").Value = 52413.855
$ws.Range("N39").Value = -53001.855
$ws.Range("H118").Value = 1065.3334
$ws.Range("I118").Value = 1065.3334
$ws.Range("K118").Value = 3196.0002
$ws.Range("M118").Value = -1953.0002
$ws.Range("H124").Value = 4800
$ws.Range("I124").Value = 4800
$ws.Range("K124").Value = 14400
$ws.Range("M124").Value = -9490
$ws.Range("H139").Value = 2426.25
$ws.Range("I139").Value = 2486.3076
$ws.Range("J139").Value = 2166
$ws.Range("K139").Value = 7458.9228
$ws.Range("L139").Value = 6498
$ws.Range("M139").Value = -2318.9228
$ws.Range("N139").Value = -16778

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5808
$ws.Range("I126").Value = 10706.5
$ws.Range("J126").Value = 2542.3333
$ws.Range("K126").Value = 32119.5
$ws.Range("L126").Value = 7626.999899999999
$ws.Range("M126").Value = -29649.5
$ws.Range("N126").Value = -12566.9999
$ws.Range("H132").Value = 5118.45
$ws.Range("J132").Value = 6342.7896
$ws.Range("L132").Value = 19028.3688
$ws.Range("N132").Value = -24088.3688

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8293
$ws.Range("I68").Value = 5951.5
$ws.Range("J68").Value = 12976
$ws.Range("K68").Value = 5951.5
$ws.Range("L68").Value = 12976
$ws.Range("M68").Value = -5202.5
$ws.Range("N68").Value = -14474
$ws.Range("H71").Value = 8293
$ws.Range("I71").Value = 5951.5
$ws.Range("J71").Value = 12976
$ws.Range("K71").Value = 29757.5
$ws.Range("L71").Value = 64880
$ws.Range("M71").Value = -26013.5
$ws.Range("N71").Value = -72368
$ws.Range("H132").Value = 3845.4783
$ws.Range("I132").Value = 2608.8333
$ws.Range("K132").Value = 7826.499899999999
$ws.Range("M132").Value = -5296.499899999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1905.5714
$ws.Range("I81").Value = 2059.75
$ws.Range("K81").Value = 4119.5
$ws.Range("M81").Value = -3058.5
$ws.Range("H84").Value = 1905.5714
$ws.Range("I84").Value = 2059.75
$ws.Range("K84").Value = 20597.5
$ws.Range("M84").Value = -15293.5
$ws.Range("H96").Value = 6015.533
$ws.Range("I96").Value = 2221.5454
$ws.Range("J96").Value = 16449
$ws.Range("K96").Value = 2221.5454
$ws.Range("L96").Value = 16449
$ws.Range("M96").Value = -848.5454
$ws.Range("N96").Value = -19195
$ws.Range("H136").Value = 3952.05
$ws.Range("I136").Value = 3331.9412
$ws.Range("K136").Value = 9995.8236
$ws.Range("M136").Value = -7445.8236
